# Add 2022-Q4 data
# 1) Insert a new "2022-Q4" sheet (copied from the "2022-Q2" template sheet so it
#    keeps identical layout/formatting), positioned right after "总计" and before
#    "2022-Q2", and fill it with the new quarter's fund holdings.
# 2) Insert a corresponding summary row into the "总计" sheet.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# Duplicate the "2022-Q2" sheet (same column headers/styles/page setup) and
# place the copy immediately before it, then rename + re-populate it.
$q2.Copy($q2)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4.Range("B2").Value = "001375"
$q4.Range("C2").Value = "金元顺安优质精选灵活配置混合C"
$q4.Range("D2:G2").NumberFormat = "@"
$q4.Range("D2").Value = "0.63"
$q4.Range("E2").Value = "68.79"
$q4.Range("F2").Value = "0.76"
$q4.Range("G2").Value = "0.0048"
$q4.Range("H2").Value = 8

$q4.Range("B3").Value = "620007"
$q4.Range("C3").Value = "金元顺安优质精选灵活配置混合A"
$q4.Range("D3:G3").NumberFormat = "@"
$q4.Range("D3").Value = "0.06"
$q4.Range("E3").Value = "68.79"
$q4.Range("F3").Value = "0.76"
$q4.Range("G3").Value = "0.0005"
$q4.Range("H3").Value = 8

# Insert the new "2022-Q4" row at the top of the "总计" table (row 2), pushing
# the existing quarters down by one row each.
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 9
$summary.Range("D6").Value = 2.19

$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.08

$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.06

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.05

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.01

$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# Restore the originally-selected tab ("2021-Q1") so the tabSelected flag ends
# up back where it started instead of on the freshly duplicated sheet.
$wb.Worksheets.Item("2021-Q1").Activate()
